# Append a new job listing row (row 73) to the Dice jobs list worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A73").Value = "Golang Developer"
$ws.Range("B73").Value = "https://www.dice.com/job-detail/8c9d94ac-871b-4972-b401-a43426b7c1d8"
$ws.Range("C73").Value = "Fremont, California"
$ws.Range("D73").Value = "Full-time, Third Party, Contract"
$ws.Range("E73").Value = "$65 - $75 per hour"
$ws.Range("F73").Value = "STAND 8"
